# "Add files via upload" — rename the "Datei"/"Fragen" worksheets to their
# English equivalents ("Working sheet"/"Questions"), make "Working sheet"
# the active tab (instead of "Dashboard"), and repoint the Dashboard charts'
# series formulas (which still hard-code the old "Fragen" sheet name) at the
# renamed "Questions" sheet.

$wb = $excel.ActiveWorkbook

$wsWorking = $wb.Worksheets.Item("Datei")
$wsQuestions = $wb.Worksheets.Item("Fragen")
$wsDashboard = $wb.Worksheets.Item("Dashboard")

# Renaming the sheets also updates the defined names (ds_salaries,
# the hidden _FilterDatabase name, etc.) that point at them by name.
$wsWorking.Name = "Working sheet"
$wsQuestions.Name = "Questions"

# Repoint every chart series on the Dashboard away from the old "Fragen"
# sheet name (sheet renames don't retroactively rewrite chart series
# formulas) at the newly renamed "Questions" sheet.
$charts = $wsDashboard.ChartObjects()
for ($i = 1; $i -le $charts.Count; $i++) {
    $chart = $charts.Item($i).Chart
    $series = $chart.SeriesCollection()
    for ($j = 1; $j -le $series.Count; $j++) {
        $s = $series.Item($j)
        $s.Formula = $s.Formula.Replace("Fragen!", "Questions!")
    }
}

# "Working sheet" (formerly "Datei") becomes the active/selected tab,
# replacing "Dashboard".
$wsWorking.Activate()
$wsWorking.Select()
